$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 18 with mapping data
$ws.Range("A18").Value = "id"
$ws.Range("B18").Value = "field_ddh_harvest_sys_id"
$ws.Range("C18").Value = $false

# Update the active selection to match the new state
$ws.Range("C19").Select()
